# Correção nos dados e início da análise PNAD 2009
#
# The sheet had two "section header" rows that carried only a label in
# column A and no data values:
#   row 5 -> "situação do domicílio"
#   row 8 -> "grandes regiões e unidades da federação"
# These rows are removed entirely (their row is deleted, shifting all
# following rows up). The row-2 header labels that become orphaned
# ("unnamed: 1_level_1" / "unnamed: 5_level_1") are replaced with "total"
# to match the corrected header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the row 5 index is not
# invalidated by the deletion of row 8.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Fix up the second header row so the previously "unnamed" placeholder
# headers read "total" like the rest of the corrected header.
$ws.Range("B2").Value = "total"
$ws.Range("C2").Value = "total"
$ws.Range("F2").Value = "total"
